# Insert a new row above row 229. This shifts the existing data rows
# 229-329 down to 230-330 (including their formatting), matching the
# new dimension A1:T330 described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(229).Insert()

# Populate the newly inserted row 229 with a full data record. Most of
# the descriptive columns (market/product metadata) are identical to
# the record that used to occupy row 229 (now shifted to row 230), so
# we copy those straight across; the measurement columns get the new
# values from the edit.
$ws.Cells.Item(229, 1).Value = 9
$ws.Cells.Item(229, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(229, 3).Value = "Metropolitana"
$ws.Cells.Item(229, 4).Value = 45202
$ws.Cells.Item(229, 5).Value = 13
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100101
$ws.Cells.Item(229, 8).Value = "Berries"
$ws.Cells.Item(229, 9).Value = 100101001
$ws.Cells.Item(229, 10).Value = "Arándano (blue)"
$ws.Cells.Item(229, 11).Value = "Sin especificar"
$ws.Cells.Item(229, 12).Value = "Primera"
$ws.Cells.Item(229, 13).Value = 70
$ws.Cells.Item(229, 14).Value = 12000
$ws.Cells.Item(229, 15).Value = 12000
$ws.Cells.Item(229, 16).Value = 12000
$ws.Cells.Item(229, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(229, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(229, 19).Value = 6000
$ws.Cells.Item(229, 20).Value = 2
